# Auto-generated edit script: update '想去人数' (column F) values
# across the four worksheets, per the commit's regenerated data dump.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 2988
$ws.Range("F9").Value = 56
$ws.Range("F13").Value = 2668
$ws.Range("F15").Value = 1536
$ws.Range("F16").Value = 7104
$ws.Range("F18").Value = 7247
$ws.Range("F21").Value = 5509
$ws.Range("F22").Value = 3119
$ws.Range("F23").Value = 3490
$ws.Range("F24").Value = 1
$ws.Range("F26").Value = 187
$ws.Range("F27").Value = 1901
$ws.Range("F29").Value = 305
$ws.Range("F30").Value = 880
$ws.Range("F34").Value = 2434
$ws.Range("F35").Value = 1214
$ws.Range("F36").Value = 2754
$ws.Range("F37").Value = 35
$ws.Range("F43").Value = 481

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 218
$ws.Range("F12").Value = 215

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 63

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 2988
$ws.Range("F8").Value = 56
$ws.Range("F12").Value = 63
$ws.Range("F13").Value = 2668
$ws.Range("F14").Value = 1536
$ws.Range("F15").Value = 218
$ws.Range("F18").Value = 7104
$ws.Range("F20").Value = 7247
$ws.Range("F22").Value = 5509
$ws.Range("F23").Value = 3119
$ws.Range("F24").Value = 215
$ws.Range("F25").Value = 3490
$ws.Range("F29").Value = 1901
$ws.Range("F32").Value = 305
$ws.Range("F33").Value = 880
$ws.Range("F37").Value = 2434
$ws.Range("F38").Value = 1214
$ws.Range("F40").Value = 2754
$ws.Range("F41").Value = 35
$ws.Range("F48").Value = 481
